# edit.ps1 — applies the "nap vegi doksi kovetes" change set to 4_select.docx
#
# Strategy: OOXML-exact edits via Range.InsertXML using a full
# pkg:package wrapper (the same mechanism real Word uses for
# "paste/insert as XML"). For every paragraph whose content changes we
# replace the paragraph's whole Range (Start..End, which includes the
# paragraph mark) with a freshly authored <w:p> so the run/proofErr
# layout matches the target exactly. New paragraphs are first created
# with Range.InsertParagraphAfter() (so a real paragraph boundary
# exists) and then their (now empty, bounded) Range is filled the same
# way.
#
# Processed highest paragraph index -> lowest so earlier indices stay
# valid for later ops.

$d = $word.ActiveDocument

$pkgHeader = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgFooter = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

function Replace-ParagraphXml($para, [string]$innerBodyXml) {
    $r = $para.Range
    $xml = $pkgHeader + $innerBodyXml + $pkgFooter
    $r.InsertXML($xml)
}

function Insert-ParagraphAfterXml($para, [string]$innerBodyXml) {
    $r = $para.Range
    $r.InsertParagraphAfter()
    # Re-fetch: the new (still empty) paragraph is now immediately
    # after $para.
    $newIndex = $para.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    Replace-ParagraphXml $newPara $innerBodyXml
}

# ---------------------------------------------------------------------
# 1) Delete the seven trailing paragraphs between the numId=6 list and
#    the closing "Minden lekerdezeshez..." paragraph:
#      "2 nezettablat hasznalo lekerdezes" (bold heading, holds the old
#      _GoBack bookmark), the two numId=5 items, "2 osszetett
#      delete/update", the two numId=9 items, and the trailing blank
#      <w:p/>.
# ---------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(20)
$pEnd = $d.Paragraphs.Item(26)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------
# 2) numId=6 item 5: " (+IN, ANY, ALL)" -> "Lakó_Legnagyobb_házak (IN)"
# ---------------------------------------------------------------------
$p19 = $d.Paragraphs.Item(19)
$inner19 = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr></w:pPr>" +
           "<w:r><w:t>Lakó_Legnagyobb_házak</w:t></w:r>" +
           "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
           "<w:r><w:t>(IN</w:t></w:r>" +
           "<w:r><w:t>)</w:t></w:r>" +
           "</w:p>"
Replace-ParagraphXml $p19 $inner19

# ---------------------------------------------------------------------
# 3) numId=6 item 4: "." -> _GoBack bookmark + "Lakó_legidősebb_Allekerdezessel (ANY)"
# ---------------------------------------------------------------------
$p18 = $d.Paragraphs.Item(18)
$inner18 = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr></w:pPr>" +
           "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
           "<w:r><w:t>Lakó_legidősebb_</w:t></w:r>" +
           "<w:proofErr w:type='spellStart'/><w:r><w:t>Allekerdezessel</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
           "<w:r><w:t xml:space='preserve'> (ANY)</w:t></w:r>" +
           "</w:p>"
Replace-ParagraphXml $p18 $inner18

# ---------------------------------------------------------------------
# 4) numId=6 item 3: "." -> "Lakó_Legfiatalabb_Allekerdezessel (ALL)"
# ---------------------------------------------------------------------
$p17 = $d.Paragraphs.Item(17)
$inner17 = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr></w:pPr>" +
           "<w:r><w:t>Lakó_Legfiatalabb_</w:t></w:r>" +
           "<w:proofErr w:type='spellStart'/><w:r><w:t>Allekerdezessel</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
           "<w:r><w:t xml:space='preserve'> (</w:t></w:r>" +
           "<w:r><w:t>ALL</w:t></w:r>" +
           "<w:r><w:t>)</w:t></w:r>" +
           "</w:p>"
Replace-ParagraphXml $p17 $inner17

# ---------------------------------------------------------------------
# 5) numId=6 item 2: "??lakó_lista_kezdőbetűk" -> "Lakó_legidősebb_maxxal"
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$inner16 = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr></w:pPr>" +
           "<w:r><w:t>Lakó_legidősebb_</w:t></w:r>" +
           "<w:proofErr w:type='spellStart'/><w:r><w:t>maxxal</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
           "</w:p>"
Replace-ParagraphXml $p16 $inner16

# ---------------------------------------------------------------------
# 6) numId=6 item 1: "??az 5 legsikeresebb csapat" -> "Lakó_legfiatalabb_minnel"
# ---------------------------------------------------------------------
$p15 = $d.Paragraphs.Item(15)
$inner15 = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr></w:pPr>" +
           "<w:r><w:t>Lakó_legfiatalabb_</w:t></w:r>" +
           "<w:proofErr w:type='spellStart'/><w:r><w:t>minnel</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
           "</w:p>"
Replace-ParagraphXml $p15 $inner15

# ---------------------------------------------------------------------
# 7) "5 allekerdezest hasznalo lekerdezes" heading: append trailing bold
#    space, then add a new paragraph "TODO: Formázott kiírás" (TODO in
#    red) right after it.
# ---------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$inner14 = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:ind w:left='0'/><w:rPr><w:b/></w:rPr></w:pPr>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>5 </w:t></w:r>" +
           "<w:proofErr w:type='spellStart'/><w:r><w:t>allekérdezést</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
           "<w:r><w:t xml:space='preserve'> használó </w:t></w:r>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t>lekérdezés</w:t></w:r>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
           "</w:p>"
Replace-ParagraphXml $p14 $inner14

$p14 = $d.Paragraphs.Item(14)
$innerTodo = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:ind w:left='0'/></w:pPr>" +
             "<w:r><w:rPr><w:color w:val='FF0000'/></w:rPr><w:t>TODO</w:t></w:r>" +
             "<w:r><w:t>: Formázott kiírás</w:t></w:r>" +
             "</w:p>"
Insert-ParagraphAfterXml $p14 $innerTodo

# ---------------------------------------------------------------------
# 8) "1: ház_lakó_lista" -> append " (inline nézettel)"
# ---------------------------------------------------------------------
$p10 = $d.Paragraphs.Item(10)
$inner10 = "<w:p><w:pPr><w:pStyle w:val='Listaszerbekezds'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
           "<w:r><w:t>1:</w:t></w:r>" +
           "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
           "<w:proofErr w:type='gramStart'/><w:r><w:t>ház_lakó</w:t></w:r><w:proofErr w:type='gramEnd'/>" +
           "<w:r><w:t>_lista</w:t></w:r>" +
           "<w:r><w:t xml:space='preserve'> (</w:t></w:r>" +
           "<w:proofErr w:type='spellStart'/><w:r><w:t>inline</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
           "<w:r><w:t xml:space='preserve'> nézettel)</w:t></w:r>" +
           "</w:p>"
Replace-ParagraphXml $p10 $inner10

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
